$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 250001010
$ws.Range("I40").Value = 2001
$ws.Range("K40").Value = 2001
$ws.Range("M40").Value = -1826

$ws.Range("H100").Value = 6857.5557
$ws.Range("I100").Value = 5713.125
$ws.Range("J100").Value = 7773.1
$ws.Range("K100").Value = 5713.125
$ws.Range("L100").Value = 7773.1
$ws.Range("M100").Value = -5172.125
$ws.Range("N100").Value = -8855.1

$ws.Range("H125").Value = 1381.4117
$ws.Range("I125").Value = 565.3333
$ws.Range("J125").Value = 1556.2858
$ws.Range("K125").Value = 5087.9997
$ws.Range("L125").Value = 14006.5722
$ws.Range("M125").Value = -2627.9997
$ws.Range("N125").Value = -18926.5722

$ws.Range("H138").Value = 4835.65
$ws.Range("J138").Value = 5603.1963
$ws.Range("L138").Value = 16809.5889
$ws.Range("N138").Value = -27089.5889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8370.091
$ws.Range("I32").Value = 8053.7764
$ws.Range("J32").Value = 17332.334
$ws.Range("K32").Value = 8053.7764
$ws.Range("L32").Value = 17332.334
$ws.Range("M32").Value = -7766.7764
$ws.Range("N32").Value = -17906.334

$ws.Range("H102").Value = 3084.3572
$ws.Range("I102").Value = 2667.7693
$ws.Range("J102").Value = 8500
$ws.Range("K102").Value = 2667.7693
$ws.Range("L102").Value = 8500
$ws.Range("M102").Value = -1045.7693
$ws.Range("N102").Value = -11744

$ws.Range("H122").Value = 5453.3335
$ws.Range("I122").Value = 4149.154
$ws.Range("J122").Value = 7572.625
$ws.Range("K122").Value = 12447.462
$ws.Range("L122").Value = 22717.875
$ws.Range("M122").Value = -9997.462000000001
$ws.Range("N122").Value = -27617.875

$ws.Range("H132").Value = 2609.861
$ws.Range("I132").Value = 2906.1724
$ws.Range("J132").Value = 1382.2858
$ws.Range("K132").Value = 8718.5172
$ws.Range("L132").Value = 4146.857400000001
$ws.Range("M132").Value = -6188.5172
$ws.Range("N132").Value = -9206.857400000001

$ws.Range("H138").Value = 47885.8
$ws.Range("J138").Value = 99429
$ws.Range("L138").Value = 99429
$ws.Range("N138").Value = -109709

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1225.3334
$ws.Range("I22").Value = 1256.125
$ws.Range("J22").Value = 979
$ws.Range("K22").Value = 1256.125
$ws.Range("L22").Value = 979
$ws.Range("M22").Value = -1083.125
$ws.Range("N22").Value = -1325

$ws.Range("H105").Value = 5730.6924
$ws.Range("I105").Value = 4983.5293
$ws.Range("J105").Value = 7142
$ws.Range("K105").Value = 4983.5293
$ws.Range("L105").Value = 7142
$ws.Range("M105").Value = -3236.5293
$ws.Range("N105").Value = -10636

$ws.Range("H134").Value = 2338.7805
$ws.Range("I134").Value = 1894.7742
$ws.Range("K134").Value = 5684.3226
$ws.Range("M134").Value = -3149.3226

$ws.Range("H141").Value = 192086.9
$ws.Range("I141").Value = 188490
$ws.Range("J141").Value = 192446.6
$ws.Range("K141").Value = 188490
$ws.Range("L141").Value = 192446.6
$ws.Range("M141").Value = -183310
$ws.Range("N141").Value = -202806.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 741.1053000000001
$ws.Range("I16").Value = 592.2353000000001
$ws.Range("K16").Value = 592.2353000000001
$ws.Range("M16").Value = -305.2353000000001

$ws.Range("H51").Value = 15224.75
$ws.Range("I51").Value = 15224.75
$ws.Range("K51").Value = 15224.75
$ws.Range("M51").Value = -14488.75

$ws.Range("H61").Value = 15224.75
$ws.Range("I61").Value = 15224.75
$ws.Range("K61").Value = 15224.75
$ws.Range("M61").Value = -14876.75

$ws.Range("H99").Value = 19203.477
$ws.Range("I99").Value = 13332.333
$ws.Range("J99").Value = 21551.934
$ws.Range("K99").Value = 13332.333
$ws.Range("L99").Value = 21551.934
$ws.Range("M99").Value = -11834.333
$ws.Range("N99").Value = -24547.934

$ws.Range("H105").Value = 4764.1055
$ws.Range("I105").Value = 965.7143
$ws.Range("J105").Value = 15399.6
$ws.Range("K105").Value = 965.7143
$ws.Range("L105").Value = 15399.6
$ws.Range("M105").Value = 781.2857
$ws.Range("N105").Value = -18893.6

$ws.Range("H113").Value = 741.1053000000001
$ws.Range("I113").Value = 592.2353000000001
$ws.Range("K113").Value = 592.2353000000001
$ws.Range("M113").Value = 1577.7647

$ws.Range("H122").Value = 5445.9
$ws.Range("I122").Value = 3522.2
$ws.Range("K122").Value = 10566.6
$ws.Range("M122").Value = -8116.599999999999

$ws.Range("H126").Value = 19203.477
$ws.Range("I126").Value = 13332.333
$ws.Range("J126").Value = 21551.934
$ws.Range("K126").Value = 39996.999
$ws.Range("L126").Value = 64655.802
$ws.Range("M126").Value = -37526.999
$ws.Range("N126").Value = -69595.802

$ws.Range("H132").Value = 2649.5715
$ws.Range("I132").Value = 2419.4
$ws.Range("K132").Value = 7258.200000000001
$ws.Range("M132").Value = -4728.200000000001

$ws.Range("H141").Value = 445997.1
$ws.Range("J141").Value = 512996.44
$ws.Range("L141").Value = 512996.44
$ws.Range("N141").Value = -523356.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2926.5098
$ws.Range("I131").Value = 2124.5
$ws.Range("J131").Value = 3075.721
$ws.Range("K131").Value = 6373.5
$ws.Range("L131").Value = 9227.163
$ws.Range("M131").Value = -1333.5
$ws.Range("N131").Value = -19307.163

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 181.76923
$ws.Range("I2").Value = 123
$ws.Range("J2").Value = 314
$ws.Range("K2").Value = 123
$ws.Range("L2").Value = 314
$ws.Range("M2").Value = -10
$ws.Range("N2").Value = -540

$ws.Range("H3").Value = 16749.25
$ws.Range("J3").Value = 16749.25
$ws.Range("L3").Value = 16749.25
$ws.Range("N3").Value = -16981.25

$ws.Range("H70").Value = 3500
$ws.Range("I70").Value = 2000
$ws.Range("K70").Value = 2000
$ws.Range("M70").Value = -1730

$ws.Range("H73").Value = 3500
$ws.Range("I73").Value = 2000
$ws.Range("K73").Value = 2000
$ws.Range("M73").Value = -1064

$ws.Range("H102").Value = 2893.1052
$ws.Range("J102").Value = 4999.25
$ws.Range("L102").Value = 4999.25
$ws.Range("N102").Value = -8243.25

$ws.Range("H122").Value = 3949.4443
$ws.Range("I122").Value = 4169.3105
$ws.Range("J122").Value = 3038.5715
$ws.Range("K122").Value = 12507.9315
$ws.Range("L122").Value = 9115.7145
$ws.Range("M122").Value = -10057.9315
$ws.Range("N122").Value = -14015.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9232
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H22").Value = 4349.5
$ws.Range("I22").Value = 4349.5
$ws.Range("K22").Value = 4349.5
$ws.Range("M22").Value = -4054.5

$ws.Range("H27").Value = 4349.5
$ws.Range("I27").Value = 4349.5
$ws.Range("K27").Value = 4349.5
$ws.Range("M27").Value = -4242.5

$ws.Range("H35").Value = 1231
$ws.Range("I35").Value = 1231
$ws.Range("K35").Value = 1231
$ws.Range("M35").Value = -895

$ws.Range("H40").Value = 5916.4814
$ws.Range("I40").Value = 5330.7144
$ws.Range("J40").Value = 7966.6665
$ws.Range("K40").Value = 5330.7144
$ws.Range("L40").Value = 7966.6665
$ws.Range("M40").Value = -5194.7144
$ws.Range("N40").Value = -8238.666499999999

$ws.Range("H82").Value = 3833.261
$ws.Range("I82").Value = 984.2857
$ws.Range("J82").Value = 8265
$ws.Range("K82").Value = 984.2857
$ws.Range("L82").Value = 8265
$ws.Range("M82").Value = -623.2857
$ws.Range("N82").Value = -8987

$ws.Range("H85").Value = 3833.261
$ws.Range("I85").Value = 984.2857
$ws.Range("J85").Value = 8265
$ws.Range("K85").Value = 984.2857
$ws.Range("L85").Value = 8265
$ws.Range("M85").Value = 263.7143
$ws.Range("N85").Value = -10761

$ws.Range("H93").Value = 1453.44
$ws.Range("I93").Value = 417.8
$ws.Range("J93").Value = 3006.9
$ws.Range("K93").Value = 417.8
$ws.Range("L93").Value = 3006.9
$ws.Range("M93").Value = 830.2
$ws.Range("N93").Value = -5502.9

$ws.Range("H126").Value = 9232
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 55164
$ws.Range("J102").Value = 55164
$ws.Range("L102").Value = 55164
$ws.Range("N102").Value = -61654

$ws.Range("H109").Value = 95000
$ws.Range("J109").Value = 95000
$ws.Range("L109").Value = 95000
$ws.Range("N109").Value = -97774

$ws.Range("H122").Value = 2283.125
$ws.Range("I122").Value = 1842
$ws.Range("J122").Value = 3959.4
$ws.Range("K122").Value = 5526
$ws.Range("L122").Value = 11878.2
$ws.Range("M122").Value = -3076
$ws.Range("N122").Value = -16778.2

$ws.Range("H136").Value = 4035.1924
$ws.Range("I136").Value = 4394.5947
$ws.Range("J136").Value = 3148.6667
$ws.Range("K136").Value = 13183.7841
$ws.Range("L136").Value = 9446.000100000001
$ws.Range("M136").Value = -10633.7841
$ws.Range("N136").Value = -14546.0001
